$d = $word.ActiveDocument

# Each of these ListBullet "answer" paragraphs gets a "Design: " prefix.
# We locate each paragraph via Find (no replacement text, so Word's
# smart-quote autocorrect never touches the existing straight quotes)
# and then prepend the prefix with InsertBefore.
$targets = @(
    'Because "cousin ship" all design material was not added into Kronodoc under NB518 because it was assumed that same persons continue from NB516-517 and they are already familiar with material',
    'Air balancing report for future USPH inspections is missing from Kronodoc.',
    '"Real" delivery control still missing. Who is in charge of TK extra prices (mainly work); approving, doing TMR and who will do final PO. ',
    'Heat load calculation "as built" version not useful; updated info in AC plans and balancing',
    'Change of system responsible person was not informed to all parties and questions were presented to "old one". Valid contact info was not available in Kronodoc. TK organisation was total mess, too many persons involved and roles unknown.',
    'Picking without part list was reduced compared to NB516-517 which was good thing.',
    'TRI / final inspection of AC rooms might be useful; has been earlier, but it was agreed during NB516 that not needed. How we confirm that all background installation done and inspected and also AC commissiong/adjustment is done before ceiling is closed ?'
)

foreach ($target in $targets) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($target, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        $range.InsertBefore("Design: ")
    }
}
